# Insert a new column before column C (everything from C..AI shifts to D..AJ),
# matching the width of the column to its left (column B), then fill in the
# new header/value cells and leave the selection on C3, as captured in the
# target commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("C").Insert()
$ws.Columns("C").ColumnWidth = 14.29

$ws.Range("C1").Value = "Project Number 2"
$ws.Range("C2").Value = "abcd"

$ws.Range("C3").Select()
